# Applies the "Atualização de bases das ligas, do dia: 28-06-2024 às 19:47" edit.
#
# The underlying change re-associates the odds/result data (columns B, F:AD)
# of a handful of fixture rows with the correct away-team / match-id, while
# each row keeps its position (column A, the positional id) and its home
# team (column E) fixed. Concretely, the match-data (everything except
# column A) is permuted among these row groups:
#   rows 10,11,12      -> rotate  (10<-12, 11<-10, 12<-11)
#   rows 25,26,27,28   -> reverse (25<-28, 26<-27, 27<-26, 28<-25)
#   rows 193,194       -> swap    (193<-194, 194<-193)
#
# We capture each row's B:AD values first (so later writes don't clobber
# data we still need to read), then write them back in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B" + $row + ":AD" + $row).Value2
}

# --- Group 1: rows 10, 11, 12 (3-way rotation) ---
$row10 = Get-RowData 10
$row11 = Get-RowData 11
$row12 = Get-RowData 12

$ws.Range("B10:AD10").Value2 = $row12
$ws.Range("B11:AD11").Value2 = $row10
$ws.Range("B12:AD12").Value2 = $row11

# --- Group 2: rows 25, 26, 27, 28 (reverse order) ---
$row25 = Get-RowData 25
$row26 = Get-RowData 26
$row27 = Get-RowData 27
$row28 = Get-RowData 28

$ws.Range("B25:AD25").Value2 = $row28
$ws.Range("B26:AD26").Value2 = $row27
$ws.Range("B27:AD27").Value2 = $row26
$ws.Range("B28:AD28").Value2 = $row25

# --- Group 3: rows 193, 194 (swap) ---
$row193 = Get-RowData 193
$row194 = Get-RowData 194

$ws.Range("B193:AD193").Value2 = $row194
$ws.Range("B194:AD194").Value2 = $row193

Write-Output "done"
